$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.316.99'
$ws.Range("E2").Value = '  +0.74%  '

# Row 3
$ws.Range("D3").Value = '3.832.32'
$ws.Range("E3").Value = '  +0.86%  '

# Row 4
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '450.22'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +6.91%  '

# Row 6
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '146.85'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +13.70%  '

# Row 7
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.624'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +3.38%  '

# Row 8
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("E9").Value = '  +3.12%  '

# Row 10
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.155'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -4.44%  '

# Row 11
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.0000322'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -8.68%  '

# Row 12
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '43.84'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +9.37%  '

# Row 13
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '10.38'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +2.78%  '

# Row 14
$ws.Range("D14").Value = '4.442.35'
$ws.Range("E14").Value = '  +1.33%  '

# Row 15
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '14.86'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -5.68%  '

# Row 16
$ws.Range("D16").Value = '3.813.86'
$ws.Range("E16").Value = '  +0.49%  '

# Row 17
$ws.Range("E17").Value = '  -0.38%  '

# Row 18
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '20.06'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +2.98%  '

# Row 19
$ws.Range("E19").Value = '  +7.97%  '

# Row 20
$ws.Range("D20").Value = '67.314.05'
$ws.Range("E20").Value = '  +0.69%  '

# Row 21
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '429.52'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +5.58%  '

# Row 22
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '14.87'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +4.26%  '

# Row 23
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '3.26'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +8.52%  '

# Row 24
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '86.52'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +3.33%  '

# Row 25
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '3.47'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +8.74%  '

# Row 26
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '37.31'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.65%  '

# Row 27
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '10.05'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +19.60%  '

# Row 28
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '9.79'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +4.16%  '

# Row 29
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '5.47'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.99%  '

# Row 30
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '735.03'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +2.81%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.135'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +12.16%  '

# Row 32
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '13.80'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +11.32%  '

# Row 33
$ws.Range("E33").Value = '  -0.88%  '

# Row 34
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '43.51'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +13.06%  '

# Row 35
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.159'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +5.41%  '

# Row 36
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '57.33'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +4.44%  '

# Row 37
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '5.59'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +12.22%  '

# Row 38
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -0.07%  '

# Row 39
$ws.Range("E39").Value = '  +5.88%  '

# Row 40
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '2.96'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +0.68%  '

# Row 41
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.349'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +17.88%  '

# Row 42
$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '0.0₃0690'
$ws.Range("E42").Value = '  -9.00%  '

# Row 43
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '2.57'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +14.21%  '

# Row 44
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.140'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +4.66%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("B46").Value = 'LidoDAOToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '3.46'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +3.78%  '

# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '3.24'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +4.84%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '2.14'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +5.52%  '

# Row 49
$ws.Range("E49").Value = '  +5.36%  '

# Row 50
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '144.97'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +0.81%  '

# Row 51
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '2.88'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +3.54%  '
